$wb = $excel.ActiveWorkbook

# --- Update "pro" sheet values (B2:B26) ---
$wsPro = $wb.Worksheets.Item("pro")
$wsPro.Range("B2").Value = 1994531.7040876676
$wsPro.Range("B3").Value = 2305342.6631757114
$wsPro.Range("B4").Value = 2538550.56879826
$wsPro.Range("B5").Value = 2609525.869466038
$wsPro.Range("B6").Value = 2032350.6249421134
$wsPro.Range("B7").Value = 1966721.5273727702
$wsPro.Range("B8").Value = 2050716.812241549
$wsPro.Range("B9").Value = 1734904.396826341
$wsPro.Range("B10").Value = 1674510.7248700035
$wsPro.Range("B11").Value = 1545219.7559931055
$wsPro.Range("B12").Value = 1603884.7505375573
$wsPro.Range("B13").Value = 1580671.2884127605
$wsPro.Range("B14").Value = 1668422.376710105
$wsPro.Range("B15").Value = 1715299.6159520743
$wsPro.Range("B16").Value = 1845730.4816233409
$wsPro.Range("B17").Value = 1510313.3664306144
$wsPro.Range("B18").Value = 1962138.207799488
$wsPro.Range("B19").Value = 2247688.7946757437
$wsPro.Range("B20").Value = 2926633.1162040355
$wsPro.Range("B21").Value = 3535999
$wsPro.Range("B22").Value = 3427804
$wsPro.Range("B23").Value = 4007293.482478193
$wsPro.Range("B24").Value = 4237914.930456318
$wsPro.Range("B25").Value = 5245494.433406818
$wsPro.Range("B26").Value = 5533996.627244192

# --- Update "ind" sheet values (B2:B101) ---
$wsInd = $wb.Worksheets.Item("ind")
$wsInd.Range("B2").Value = 466950.9927981131
$wsInd.Range("B3").Value = 412949.8575765626
$wsInd.Range("B4").Value = 590835.9500710818
$wsInd.Range("B5").Value = 581306.3379731611
$wsInd.Range("B6").Value = 898960.0745705168
$wsInd.Range("B7").Value = 689308.6084162618
$wsInd.Range("B8").Value = 673425.921586394
$wsInd.Range("B9").Value = 765614.9377103586
$wsInd.Range("B10").Value = 997328.5841495992
$wsInd.Range("B11").Value = 1042841.594934204
$wsInd.Range("B12").Value = 1116058.177500742
$wsInd.Range("B13").Value = 1149698.2289502325
$wsInd.Range("B14").Value = 1080439.2994953992
$wsInd.Range("B15").Value = 1044820.4214900563
$wsInd.Range("B16").Value = 1102206.3916097754
$wsInd.Range("B17").Value = 1670129.6131394082
$wsInd.Range("B18").Value = 1080439.299495399
$wsInd.Range("B19").Value = 1335707.9252003562
$wsInd.Range("B20").Value = 809340.0613436232
$wsInd.Range("B21").Value = 718314.0397744138
$wsInd.Range("B22").Value = 653012.7634312852
$wsInd.Range("B23").Value = 850895.419016523
$wsInd.Range("B24").Value = 844958.9393489659
$wsInd.Range("B25").Value = 1003265.0638171561
$wsInd.Range("B26").Value = 912239.0422479468
$wsInd.Range("B27").Value = 1155634.7086177894
$wsInd.Range("B28").Value = 941921.4405857326
$wsInd.Range("B29").Value = 524389.0373008806
$wsInd.Range("B30").Value = 736123.478777085
$wsInd.Range("B31").Value = 823191.8472345899
$wsInd.Range("B32").Value = 997328.5841495992
$wsInd.Range("B33").Value = 755911.7443356091
$wsInd.Range("B34").Value = 805382.4082319182
$wsInd.Range("B35").Value = 1007222.7169288611
$wsInd.Range("B36").Value = 841001.2862372613
$wsInd.Range("B37").Value = 627288.0182052043
$wsInd.Range("B38").Value = 832510.4650989077
$wsInd.Range("B39").Value = 871399.3622575073
$wsInd.Range("B40").Value = 675514.5469401171
$wsInd.Range("B41").Value = 648148.2859766582
$wsInd.Range("B42").Value = 882556.6439101612
$wsInd.Range("B43").Value = 978133.9665578315
$wsInd.Range("B44").Value = 813297.714455328
$wsInd.Range("B45").Value = 919560.7005046009
$wsInd.Range("B46").Value = 884860.0848999444
$wsInd.Range("B47").Value = 832007.2919771483
$wsInd.Range("B48").Value = 878643.5112432821
$wsInd.Range("B49").Value = 946027.7206193072
$wsInd.Range("B50").Value = 1087610.1664836933
$wsInd.Range("B51").Value = 1018309.5748498221
$wsInd.Range("B52").Value = 887610.2016270369
$wsInd.Range("B53").Value = 840596.6953594033
$wsInd.Range("B54").Value = 1046261.7231959478
$wsInd.Range("B55").Value = 964847.8959976218
$wsInd.Range("B56").Value = 1001721.4634724552
$wsInd.Range("B57").Value = 1088213.360243661
$wsInd.Range("B58").Value = 1328130.431841921
$wsInd.Range("B59").Value = 1211400.6755483735
$wsInd.Range("B60").Value = 1209578.172999187
$wsInd.Range("B61").Value = 1040604.0296953882
$wsInd.Range("B62").Value = 778634.3183386849
$wsInd.Range("B63").Value = 627544.6132763181
$wsInd.Range("B64").Value = 1198758.561172776
$wsInd.Range("B65").Value = 1548315.1853982068
$wsInd.Range("B66").Value = 1840731.6404052698
$wsInd.Range("B67").Value = 1654266.2165459578
$wsInd.Range("B68").Value = 1621277.812980243
$wsInd.Range("B69").Value = 1508296.8599212281
$wsInd.Range("B70").Value = 1699059.2718491217
$wsInd.Range("B71").Value = 2127371.99341605
$wsInd.Range("B72").Value = 1918676.8829231742
$wsInd.Range("B73").Value = 2527254.265868389
$wsInd.Range("B74").Value = 2083536.0372116824
$wsInd.Range("B75").Value = 2117873.8620703393
$wsInd.Range("B76").Value = 2209443.8498000987
$wsInd.Range("B77").Value = 2357623.0598190464
$wsInd.Range("B78").Value = 2525797.9540317263
$wsInd.Range("B79").Value = 3024420.0557582094
$wsInd.Range("B80").Value = 2733548.9853139445
$wsInd.Range("B81").Value = 2729904.767484191
$wsInd.Range("B82").Value = 2816164.183171855
$wsInd.Range("B83").Value = 2766385.5836006426
$wsInd.Range("B84").Value = 2563415.001496254
$wsInd.Range("B85").Value = 2530709.1002720627
$wsInd.Range("B86").Value = 2831074.884191337
$wsInd.Range("B87").Value = 3042910.164803041
$wsInd.Range("B88").Value = 2758856.0869191336
$wsInd.Range("B89").Value = 3159880.740872405
$wsInd.Range("B90").Value = 3407924.578789104
$wsInd.Range("B91").Value = 3706950.307152768
$wsInd.Range("B92").Value = 3503083.5493119294
$wsInd.Range("B93").Value = 3112619.2088785116
$wsInd.Range("B94").Value = 3394695.955137011
$wsInd.Range("B95").Value = 4239101.086167534
$wsInd.Range("B96").Value = 4046780.873013761
$wsInd.Range("B97").Value = 5297368.424773461
$wsInd.Range("B98").Value = 3954431.0747140837
$wsInd.Range("B99").Value = 3921849.3492875877
$wsInd.Range("B100").Value = 4156684.727273116
$wsInd.Range("B101").Value = 3653967.379703049

# --- Update "conso" sheet values (B2:B26) ---
$wsConso = $wb.Worksheets.Item("conso")
$wsConso.Range("B2").Value = 1306657.5890816844
$wsConso.Range("B3").Value = 1510279.947189798
$wsConso.Range("B4").Value = 1663059.9640289014
$wsConso.Range("B5").Value = 1709555.3416434436
$wsConso.Range("B6").Value = 1331431.414550583
$wsConso.Range("B7").Value = 1288435.715171447
$wsConso.Range("B8").Value = 1343460.1183232337
$wsConso.Range("B9").Value = 1136560.4169311868
$wsConso.Range("B10").Value = 1096996.0276962311
$wsConso.Range("B11").Value = 1012295.6163525191
$wsConso.Range("B12").Value = 1050727.3671621352
$wsConso.Range("B13").Value = 1035519.6707268794
$wsConso.Range("B14").Value = 1093006.145804898
$wsConso.Range("B15").Value = 1123715.664686847
$wsConso.Range("B16").Value = 1209161.8449723872
$wsConso.Range("B17").Value = 989423.6087723221
$wsConso.Range("B18").Value = 1285418.5368348253
$wsConso.Range("B19").Value = 1472486.3046891817
$wsConso.Range("B20").Value = 1917272.2850479844
$wsConso.Range("B21").Value = 2316476
$wsConso.Range("B22").Value = 1790216
$wsConso.Range("B23").Value = 2022744.696644935
$wsConso.Range("B24").Value = 2132026.1055200375
$wsConso.Range("B25").Value = 2638948.7875440973
$wsConso.Range("B26").Value = 2784090.9708590223

# --- Recalculate so VA (= pro - conso) picks up the new values ---
$excel.Calculate() | Out-Null

# --- Update sheet selections / active sheet (D24 on every sheet; "pro" tab active) ---
foreach ($name in @("pro","ind","VA","conso")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate() | Out-Null
    $ws.Range("D24").Select() | Out-Null
}
$wb.Worksheets.Item("pro").Activate() | Out-Null
$wb.Worksheets.Item("pro").Range("D24").Select() | Out-Null
